$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 1.018777190134747
$ws.Range("E2").Value = 1.028441947054262
$ws.Range("F2").Value = 1.035878506188935
$ws.Range("J2").Value = 1.023983701230586
$ws.Range("L2").Value = 1.031258752783232
$ws.Range("M2").Value = 1.038673849366234
$ws.Range("N2").Value = 1.012006349856506
$ws.Range("C3").Value = 1.019837964448458
$ws.Range("E3").Value = 1.029421715039732
$ws.Range("F3").Value = 1.037028821821524
$ws.Range("J3").Value = 1.024680361429499
$ws.Range("L3").Value = 1.032045993535909
$ws.Range("M3").Value = 1.039632745861786
$ws.Range("N3").Value = 1.012244501852444
$ws.Range("C4").Value = 1.020524907074266
$ws.Range("E4").Value = 1.030056559725999
$ws.Range("F4").Value = 1.037774361322412
$ws.Range("J4").Value = 1.025131167901118
$ws.Range("L4").Value = 1.032555648281803
$ws.Range("M4").Value = 1.040253847108459
$ws.Range("N4").Value = 1.01239841392067
$ws.Range("C5").Value = 1.020813829861862
$ws.Range("E5").Value = 1.030323656259277
$ws.Range("F5").Value = 1.038088075393573
$ws.Range("J5").Value = 1.025320691769214
$ws.Range("L5").Value = 1.032769968569184
$ws.Range("M5").Value = 1.040515109076831
$ws.Range("N5").Value = 1.012463073293099
$ws.Range("C6").Value = 1.020862349024004
$ws.Range("E6").Value = 1.030368515121524
$ws.Range("F6").Value = 1.038140766328743
$ws.Range("J6").Value = 1.025352513939803
$ws.Range("L6").Value = 1.032805957485931
$ws.Range("M6").Value = 1.040558984963684
$ws.Range("N6").Value = 1.012473927226143
$ws.Range("C7").Value = 1.02052876715525
$ws.Range("E7").Value = 1.030060127866673
$ws.Range("F7").Value = 1.037778552048582
$ws.Range("J7").Value = 1.025133700309516
$ws.Range("L7").Value = 1.032558511798797
$ws.Range("M7").Value = 1.040257337510414
$ws.Range("N7").Value = 1.012399278080471
$ws.Range("C8").Value = 1.01913556990954
$ws.Range("E8").Value = 1.028772884070614
$ws.Range("F8").Value = 1.036267010527894
$ws.Range("J8").Value = 1.024219136428681
$ws.Range("L8").Value = 1.031524751091472
$ws.Range("M8").Value = 1.038997782203663
$ws.Range("N8").Value = 1.012086873282085
$ws.Range("C9").Value = 1.01668478420329
$ws.Range("E9").Value = 1.026511265982619
$ws.Range("F9").Value = 1.033612736623003
$ws.Range("J9").Value = 1.022607722786081
$ws.Range("L9").Value = 1.029705111133826
$ws.Range("M9").Value = 1.036783121677486
$ws.Range("N9").Value = 1.011534943782861
$ws.Range("C10").Value = 1.015053735864303
$ws.Range("E10").Value = 1.025008019577426
$ws.Range("F10").Value = 1.03184944187669
$ws.Range("J10").Value = 1.021533564086832
$ws.Range("L10").Value = 1.028493356457054
$ws.Range("M10").Value = 1.035309935747672
$ws.Range("N10").Value = 1.011166039416976
$ws.Range("C11").Value = 1.014348135187404
$ws.Range("E11").Value = 1.024358166769749
$ws.Range("F11").Value = 1.031087390176909
$ws.Range("J11").Value = 1.021068470029094
$ws.Range("L11").Value = 1.027968971972987
$ws.Range("M11").Value = 1.034672801649342
$ws.Range("N11").Value = 1.011006076378505
$ws.Range("C12").Value = 1.014086141422113
$ws.Range("E12").Value = 1.024116942290496
$ws.Range("F12").Value = 1.0308045503853
$ws.Range("J12").Value = 1.020895716942084
$ws.Range("L12").Value = 1.027774239312679
$ws.Range("M12").Value = 1.03443625652106
$ws.Range("N12").Value = 1.010946625336262
$ws.Range("C13").Value = 1.01414233555925
$ws.Range("E13").Value = 1.024168678527527
$ws.Range("F13").Value = 1.030865210521971
$ws.Range("J13").Value = 1.020932772905002
$ws.Range("L13").Value = 1.027816007998315
$ws.Range("M13").Value = 1.034486991058778
$ws.Range("N13").Value = 1.010959379306702
$ws.Range("C14").Value = 1.014326476710061
$ws.Range("E14").Value = 1.024338223818738
$ws.Range("F14").Value = 1.03106400606181
$ws.Range("J14").Value = 1.021054190132586
$ws.Range("L14").Value = 1.027952874349895
$ws.Range("M14").Value = 1.034653246413062
$ws.Range("N14").Value = 1.011001162822286
$ws.Range("C15").Value = 1.014439945087112
$ws.Range("E15").Value = 1.024442707445325
$ws.Range("F15").Value = 1.03118651971999
$ws.Range("J15").Value = 1.021128999764434
$ws.Range("L15").Value = 1.028037208466125
$ws.Range("M15").Value = 1.034755697038808
$ws.Range("N15").Value = 1.011026902571408
$ws.Range("C16").Value = 1.015100577992575
$ws.Range("E16").Value = 1.025051170569725
$ws.Range("F16").Value = 1.031900047620276
$ws.Range("J16").Value = 1.02156443134678
$ws.Range("L16").Value = 1.028528164730228
$ws.Range("M16").Value = 1.035352236350497
$ws.Range("N16").Value = 1.011176650909738
$ws.Range("C17").Value = 1.015515150293139
$ws.Range("E17").Value = 1.025433128034535
$ws.Range("F17").Value = 1.032348017322074
$ws.Range("J17").Value = 1.021837572421634
$ws.Range("L17").Value = 1.028836212510739
$ws.Range("M17").Value = 1.03572663446968
$ws.Range("N17").Value = 1.01127052395519
$ws.Range("C18").Value = 1.015757026439381
$ws.Range("E18").Value = 1.025656019899016
$ws.Range("F18").Value = 1.032609451984647
$ws.Range("J18").Value = 1.021996893281358
$ws.Range("L18").Value = 1.029015921795614
$ws.Range("M18").Value = 1.03594508838923
$ws.Range("N18").Value = 1.011325256793188
$ws.Range("C19").Value = 1.015839510679601
$ws.Range("E19").Value = 1.025732037653505
$ws.Range("F19").Value = 1.032698618457841
$ws.Range("J19").Value = 1.02205121796653
$ws.Range("L19").Value = 1.029077203144092
$ws.Range("M19").Value = 1.036019588056826
$ws.Range("N19").Value = 1.01134391558742
$ws.Range("C20").Value = 1.015470664081779
$ws.Range("E20").Value = 1.025392137006533
$ws.Range("F20").Value = 1.032299939776394
$ws.Range("J20").Value = 1.021808266732241
$ws.Range("L20").Value = 1.028803158772577
$ws.Range("M20").Value = 1.035686457458765
$ws.Range("N20").Value = 1.011260454511378
$ws.Range("C21").Value = 1.014272249036662
$ws.Range("E21").Value = 1.024288292544409
$ws.Range("F21").Value = 1.03100545965476
$ws.Range("J21").Value = 1.021018435687277
$ws.Range("L21").Value = 1.027912569320712
$ws.Range("M21").Value = 1.034604285185607
$ws.Range("N21").Value = 1.010988859545931
$ws.Range("C22").Value = 1.013519323194119
$ws.Range("E22").Value = 1.023595186264539
$ws.Range("F22").Value = 1.030192842195708
$ws.Range("J22").Value = 1.020521858121802
$ws.Range("L22").Value = 1.027352892817892
$ws.Range("M22").Value = 1.033924544930364
$ws.Range("N22").Value = 1.010817902416904
$ws.Range("C23").Value = 1.013918410019024
$ws.Range("E23").Value = 1.023962527342229
$ws.Range("F23").Value = 1.030623505429295
$ws.Range("J23").Value = 1.020785101339317
$ws.Range("L23").Value = 1.027649562091004
$ws.Range("M23").Value = 1.034284825157512
$ws.Range("N23").Value = 1.010908548417934
$ws.Range("C24").Value = 1.015490765290502
$ws.Range("E24").Value = 1.025410658770727
$ws.Range("F24").Value = 1.032321663512391
$ws.Range("J24").Value = 1.021821508704645
$ws.Range("L24").Value = 1.028818094240573
$ws.Range("M24").Value = 1.035704611492334
$ws.Range("N24").Value = 1.011265004526705
$ws.Range("C25").Value = 1.017317874621408
$ws.Range("E25").Value = 1.027095156568682
$ws.Range("F25").Value = 1.034297833428939
$ws.Range("J25").Value = 1.023024291679433
$ws.Range("L25").Value = 1.030175296368113
$ws.Range("M25").Value = 1.037355091090548
$ws.Range("N25").Value = 1.011677799327959
